# The author renamed the "Perturbation_*" vocabulary used throughout the
# data-dictionary (shared strings / column headers) to "Disturbance_*".
# This is exactly what Excel's Find & Replace ("Replace All") across the
# whole sheet would do. Two passes are needed (match-case, whole workbook)
# so that both the Title-case tokens (Perturbation_date_day, ...,
# "Time after perturbation" the specific lowercase occurrence) map onto
# the correctly-cased replacement ("Disturbance_date_day", ...,
# "Time after disturbance").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the capitalized form: Perturbation -> Disturbance
$ws.Cells.Replace(
    "Perturbation",
    "Disturbance",
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart,
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows,
    $true,
    $false,
    $false
)

# Replace the lowercase form: perturbation -> disturbance
$ws.Cells.Replace(
    "perturbation",
    "disturbance",
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart,
    [Microsoft.Office.Interop.Excel.XlSearchOrder]::xlByRows,
    $true,
    $false,
    $false
)

# Restore the cursor/scroll position left behind in the saved file
# (active cell AJ36, viewport scrolled so AH21 is the top-left cell).
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 34
$ws.Range("AJ36").Select()
